$wb = $excel.ActiveWorkbook

# xlEdgeTop = 8, xlEdgeBottom = 9, xlEdgeRight = 10 (standard Excel Borders(Index) constants)
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10
$xlContinuous = 1

function Set-TopBottomBorder($cell) {
    # Start from the default (unstyled) state so the new xf only carries the
    # border - matches the "approach" cells which get a plain top+bottom rule
    # with no bold/center formatting inherited from the header style.
    $cell.ClearFormats()
    $cell.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
    $cell.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
}

function Set-TopBottomRightBorder($cell) {
    $cell.ClearFormats()
    $cell.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
    $cell.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
    $cell.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
}

# ---- Sheet "quality_comparison" ----
$ws1 = $wb.Worksheets.Item("quality_comparison")

Set-TopBottomBorder $ws1.Range("C1")
Set-TopBottomRightBorder $ws1.Range("D1")

$ws1.Range("C2").Value = "approach"

# ---- Sheet "computational_comparison" ----
$ws2 = $wb.Worksheets.Item("computational_comparison")

Set-TopBottomBorder $ws2.Range("C1")
Set-TopBottomRightBorder $ws2.Range("D1")
Set-TopBottomBorder $ws2.Range("F1")
Set-TopBottomRightBorder $ws2.Range("G1")

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# G5 held an empty inline string placeholder; drop it back to a truly blank cell.
$ws2.Range("G5").ClearContents()
